$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F28").Value = 40
$ws.Range("G28").Value = 4098
$ws.Range("F35").Value = 202
$ws.Range("G35").Value = 5692.36
$ws.Range("B40").Value = 75023.28999999999
$ws.Range("F63").Value = 49
$ws.Range("G63").Value = 1093.68
$ws.Range("B73").Value = 265980.08
$ws.Range("F114").Value = 20
$ws.Range("G114").Value = 18338.8
$ws.Range("B115").Value = 25816.43
$ws.Range("B149").Value = 65258
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 64287.16
$ws.Range("B150").Value = 64196
$ws.Range("F150").Value = 1
$ws.Range("G150").Value = 32143.58
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("B162").Value = 0
$ws.Range("F260").Value = 0
$ws.Range("G260").Value = 0
$ws.Range("F261").Value = 0
$ws.Range("G261").Value = 0
$ws.Range("B262").Value = 853.45
$ws.Range("F285").Value = 1940
$ws.Range("G285").Value = 35890
$ws.Range("B292").Value = 56310.17
$ws.Range("F329").Value = 48
$ws.Range("G329").Value = 5037.12
$ws.Range("F337").Value = 172
$ws.Range("G337").Value = 7306.56
$ws.Range("F342").Value = 40
$ws.Range("G342").Value = 4582.4
$ws.Range("F356").Value = 33
$ws.Range("G356").Value = 3204.3
$ws.Range("F360").Value = 198
$ws.Range("G360").Value = 9280.26
$ws.Range("B370").Value = 64985
$ws.Range("C370").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F370").Value = 13
$ws.Range("G370").Value = 1140.1
$ws.Range("B371").Value = 66196
$ws.Range("C371").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F371").Value = 28
$ws.Range("G371").Value = 2455.6
$ws.Range("F372").Value = 5
$ws.Range("G372").Value = 2116.05
$ws.Range("B373").Value = 148489.65
$ws.Range("B398").Value = 60325
$ws.Range("E398").Value = 151.57
$ws.Range("F398").Value = -102
$ws.Range("G398").Value = -12939.72
$ws.Range("B399").Value = 63560
$ws.Range("E399").Value = 134.87
$ws.Range("F399").Value = 1
$ws.Range("G399").Value = 126.86
$ws.Range("F406").Value = 176
$ws.Range("G406").Value = 30154.08
$ws.Range("B410").Value = 39943
$ws.Range("B481").Value = 58047
$ws.Range("D481").Value = 105.54
$ws.Range("E481").Value = 126.1
$ws.Range("F481").Value = 34
$ws.Range("G481").Value = 3588.36
$ws.Range("B482").Value = 47097
$ws.Range("D482").Value = 112.28
$ws.Range("E482").Value = 134.16
$ws.Range("F482").Value = 15
$ws.Range("G482").Value = 1684.2
$ws.Range("B551").Value = 53263
$ws.Range("E551").Value = 15.29
$ws.Range("F551").Value = -309
$ws.Range("G551").Value = -3958.29
$ws.Range("B552").Value = 65066
$ws.Range("E552").Value = 13.61
$ws.Range("F552").Value = 90
$ws.Range("G552").Value = 1152.9
$ws.Range("B560").Value = 45718
$ws.Range("E560").Value = 19.38
$ws.Range("F560").Value = -294
$ws.Range("G560").Value = -4768.68
$ws.Range("B561").Value = 64927
$ws.Range("E561").Value = 17.26
$ws.Range("F561").Value = 106
$ws.Range("G561").Value = 1719.32
$ws.Range("B570").Value = 65067
$ws.Range("E570").Value = 15.65
$ws.Range("F570").Value = 126
$ws.Range("G570").Value = 1855.98
$ws.Range("B571").Value = 53595
$ws.Range("E571").Value = 17.61
$ws.Range("F571").Value = -335
$ws.Range("G571").Value = -4934.55
$ws.Range("F590").Value = 0
$ws.Range("G590").Value = 0
$ws.Range("F591").Value = 0
$ws.Range("G591").Value = 0
$ws.Range("B592").Value = 0
$ws.Range("F599").Value = 0
$ws.Range("G599").Value = 0
$ws.Range("F600").Value = 0
$ws.Range("G600").Value = 0
$ws.Range("B601").Value = 322.4
$ws.Range("F636").Value = 0
$ws.Range("G636").Value = 0
$ws.Range("B637").Value = 0
$ws.Range("F639").Value = 0
$ws.Range("G639").Value = 0
$ws.Range("B644").Value = 4635.18
$ws.Range("B670").Value = 60022
$ws.Range("E670").Value = 37.22
$ws.Range("F670").Value = -113
$ws.Range("G670").Value = -3709.79
$ws.Range("B671").Value = 64830
$ws.Range("E671").Value = 34.9
$ws.Range("F671").Value = 89
$ws.Range("G671").Value = 2921.87
$ws.Range("F677").Value = 233
$ws.Range("G677").Value = 37391.84
$ws.Range("B693").Value = 167444.21
$ws.Range("F695").Value = 0
$ws.Range("G695").Value = 0
$ws.Range("B696").Value = 0
$ws.Range("F751").Value = 15
$ws.Range("G751").Value = 3766.5
$ws.Range("F754").Value = 17
$ws.Range("G754").Value = 9548.049999999999
$ws.Range("B757").Value = 55970.23
$ws.Range("F907").Value = 24
$ws.Range("G907").Value = 3114.96
$ws.Range("B913").Value = 17497.31
$ws.Range("F916").Value = 12
$ws.Range("G916").Value = 4857.24
$ws.Range("F933").Value = 7
$ws.Range("G933").Value = 1498.42
$ws.Range("B936").Value = 96846.32000000001
$ws.Range("B942").Value = 4318728.68
$ws.Range("B943").Value = 4318728.68
